# chore: changed workbook sheet and table naming
# fix: fixed converters and services affected by workbook naming changes

$wb = $excel.ActiveWorkbook

# Rename worksheets to their new, more descriptive names.
$wb.Worksheets.Item("progw_summary").Name   = "program_summary_sheet"
$wb.Worksheets.Item("progw_promoters").Name = "promoter_sheet"

# The promoters table's header label changes from the plural
# "promoters_table" to the singular "promoter_table", and is styled
# bold/Arial to match the other sheet headers.
$ws = $wb.Worksheets.Item("promoter_sheet")
$cell = $ws.Range("A1")
$cell.Value = "promoter_table"
$cell.Font.Bold = $true
$cell.Font.Name = "Arial"

Write-Output "done"
